$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Rename the header labels for the shift start/end columns
$ws.Range("I1").Value = "ShiftBegins"
$ws.Range("J1").Value = "ShiftsEnds"

# Update the active cell selection on the Events sheet
$ws.Activate()
$ws.Range("J4").Select()
